# fix(gui) step 1 and 2
# Updates the price list date and the "step 1 & 2" prices in column D
# (rows 29-37) on the GANCHO J worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Bump the list date shown in A1 (2024-01-17 -> 2024-01-18)
$ws.Range("A1").Value = 45309

# Updated prices for GANJ-50 .. GANJ-160 (column D)
$ws.Range("D29").Value = 94.935
$ws.Range("D30").Value = 101.646
$ws.Range("D31").Value = 106.442
$ws.Range("D32").Value = 111.716
$ws.Range("D33").Value = 116.993
$ws.Range("D34").Value = 123.699
$ws.Range("D35").Value = 134.252
$ws.Range("D36").Value = 151.514
$ws.Range("D37").Value = 177.405
